# Automatische test-sync: 2025-08-05 18:41:50
# Append Testmail #17 as a new row (38) to the "Logs" sheet and bump the
# "Planning / Afspraak" counter on the "Dashboard" sheet from 20 to 21.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$row = 38

$ws.Cells.Item($row, 1).Value  = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($row, 2).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 3).Value  = "Testmail #17: Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($row, 4).Value  = "Planning / Afspraak"
$ws.Cells.Item($row, 5).Value  = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($row, 6).Value  = "2025-08-05 18:41:37"
$ws.Cells.Item($row, 7).Value  = "Ja"
$ws.Cells.Item($row, 8).Value  = "Ja"
$ws.Cells.Item($row, 9).Value  = "Nee"
$ws.Cells.Item($row, 10).Value = "Nee"

# The sheet's used range (<dimension>) grows to A1:J38 automatically, but the
# conditionalFormatting sqref ranges are fixed ranges that must be extended
# explicitly to include the new row.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "37")
    $newRange = $ws.Range($col + "2:" + $col + "38")
    $fcs = $oldRange.FormatConditions
    $fcs.Item(1).ModifyAppliesToRange($newRange)
}

# Bump the "Planning / Afspraak" tally on the Dashboard sheet.
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(2, 2).Value = 21
